$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell C1 "Electrode Locations", matching the style of A1/B1
$ws.Range("A1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null
$ws.Range("C1").Value = "Electrode Locations"

# Rewrite data rows 2-27, sorted by electrode location (column letter, then row number),
# and populate the new column C with the electrode location string.
$ws.Cells.Item(2, 1).Value = "C3_bipolar_10V_5Hz.txt"
$ws.Cells.Item(2, 2).Value = 104.107638
$ws.Cells.Item(2, 3).Value = "C3"

$ws.Cells.Item(3, 1).Value = "C6_bipolar_10V_5Hz.txt"
$ws.Cells.Item(3, 2).Value = 124.69703
$ws.Cells.Item(3, 3).Value = "C6"

$ws.Cells.Item(4, 1).Value = "C8_bipolar_10V_5Hz.txt"
$ws.Cells.Item(4, 2).Value = 106.420917
$ws.Cells.Item(4, 3).Value = "C8"

$ws.Cells.Item(5, 1).Value = "C10_bipolar_10V_5Hz.txt"
$ws.Cells.Item(5, 2).Value = 92.442334
$ws.Cells.Item(5, 3).Value = "C10"

$ws.Cells.Item(6, 1).Value = "C13_bipolar_10V_5Hz.txt"
$ws.Cells.Item(6, 2).Value = 87.82990700000001
$ws.Cells.Item(6, 3).Value = "C13"

$ws.Cells.Item(7, 1).Value = "E3_bipolar_10V_5Hz.txt"
$ws.Cells.Item(7, 2).Value = 94.834337
$ws.Cells.Item(7, 3).Value = "E3"

$ws.Cells.Item(8, 1).Value = "E12_bipolar_10V_5Hz.txt"
$ws.Cells.Item(8, 2).Value = 82.462535
$ws.Cells.Item(8, 3).Value = "E12"

$ws.Cells.Item(9, 1).Value = "E13_bipolar_10V_5Hz.txt"
$ws.Cells.Item(9, 2).Value = 87.442342
$ws.Cells.Item(9, 3).Value = "E13"

$ws.Cells.Item(10, 1).Value = "E15_bipolar_10V_5Hz.txt"
$ws.Cells.Item(10, 2).Value = 77.27683399999999
$ws.Cells.Item(10, 3).Value = "E15"

$ws.Cells.Item(11, 1).Value = "G1_bipolar_10V_5Hz.txt"
$ws.Cells.Item(11, 2).Value = 70.829528
$ws.Cells.Item(11, 3).Value = "G1"

$ws.Cells.Item(12, 1).Value = "G3_bipolar_10V_5Hz.txt"
$ws.Cells.Item(12, 2).Value = 75.898152
$ws.Cells.Item(12, 3).Value = "G3"

$ws.Cells.Item(13, 1).Value = "G5_bipolar_10V_5Hz.txt"
$ws.Cells.Item(13, 2).Value = 83.29216599999999
$ws.Cells.Item(13, 3).Value = "G5"

$ws.Cells.Item(14, 1).Value = "G10_bipolar_10V_5Hz.txt"
$ws.Cells.Item(14, 2).Value = 84.039036
$ws.Cells.Item(14, 3).Value = "G10"

$ws.Cells.Item(15, 1).Value = "G15_bipolar_10V_5Hz.txt"
$ws.Cells.Item(15, 2).Value = 70.026138
$ws.Cells.Item(15, 3).Value = "G15"

$ws.Cells.Item(16, 1).Value = "H14_bipolar_10V_5Hz.txt"
$ws.Cells.Item(16, 2).Value = 78.346675
$ws.Cells.Item(16, 3).Value = "H14"

$ws.Cells.Item(17, 1).Value = "I12_bipolar_10V_5Hz.txt"
$ws.Cells.Item(17, 2).Value = 95.779027
$ws.Cells.Item(17, 3).Value = "I12"

$ws.Cells.Item(18, 1).Value = "I13_bipolar_10V_5Hz.txt"
$ws.Cells.Item(18, 2).Value = 99.111682
$ws.Cells.Item(18, 3).Value = "I13"

$ws.Cells.Item(19, 1).Value = "K1_bipolar_10V_5Hz.txt"
$ws.Cells.Item(19, 2).Value = 80.982924
$ws.Cells.Item(19, 3).Value = "K1"

$ws.Cells.Item(20, 1).Value = "K5_bipolar_10V_5Hz.txt"
$ws.Cells.Item(20, 2).Value = 86.420946
$ws.Cells.Item(20, 3).Value = "K5"

$ws.Cells.Item(21, 1).Value = "K13_bipolar_10V_5Hz.txt"
$ws.Cells.Item(21, 2).Value = 73.93812200000001
$ws.Cells.Item(21, 3).Value = "K13"

$ws.Cells.Item(22, 1).Value = "M8_bipolar_10V_5Hz.txt"
$ws.Cells.Item(22, 2).Value = 94.949395
$ws.Cells.Item(22, 3).Value = "M8"

$ws.Cells.Item(23, 1).Value = "M10_bipolar_10V_5Hz.txt"
$ws.Cells.Item(23, 2).Value = 82.83597
$ws.Cells.Item(23, 3).Value = "M10"

$ws.Cells.Item(24, 1).Value = "M12_bipolar_10V_5Hz.txt"
$ws.Cells.Item(24, 2).Value = 84.236856
$ws.Cells.Item(24, 3).Value = "M12"

$ws.Cells.Item(25, 1).Value = "O4_bipolar_10V_5Hz.txt"
$ws.Cells.Item(25, 2).Value = 97.22633999999999
$ws.Cells.Item(25, 3).Value = "O4"

$ws.Cells.Item(26, 1).Value = "O8_bipolar_10V_5Hz.txt"
$ws.Cells.Item(26, 2).Value = 122.236396
$ws.Cells.Item(26, 3).Value = "O8"

$ws.Cells.Item(27, 1).Value = "O10_bipolar_10V_5Hz.txt"
$ws.Cells.Item(27, 2).Value = 115.554936
$ws.Cells.Item(27, 3).Value = "O10"
